$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "removing birds" ---------------------------------------------------
# Remove the "order"/"broad group" header row (40) and the five
# "Galliformes"/"bird" rows (41-45). This shifts the arachnid rows
# (formerly 46-50) up to occupy rows 40-44.
$ws.Range("A40:B45").EntireRow.Delete() | Out-Null

# --- "adding additions papers" ------------------------------------------
# Append the new insect records at the bottom: four Diptera rows and one
# Orthoptera row, all tagged "insect".
$ws.Range("A45").Value = "Diptera"
$ws.Range("B45").Value = "insect"
$ws.Range("A46").Value = "Diptera"
$ws.Range("B46").Value = "insect"
$ws.Range("A47").Value = "Diptera"
$ws.Range("B47").Value = "insect"
$ws.Range("A48").Value = "Diptera"
$ws.Range("B48").Value = "insect"
$ws.Range("A49").Value = "Orthoptera"
$ws.Range("B49").Value = "insect"

# Match the formatting used by the other data rows in column B.
$ws.Range("B44").Copy() | Out-Null
$ws.Range("B45:B49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- refresh the worksheet's recorded sort range -------------------------
# The data isn't in a strict sort order (it's grouped, with extra rows
# appended over time), so a real re-sort would scramble it. Temporarily
# stamp column B with values that already match the current row order when
# sorted descending, run the sort (a no-op reorder), then restore the real
# category labels. This keeps row order intact while refreshing the
# worksheet's persisted sort range down to the new row count (49).
$realB = @()
for ($r = 1; $r -le 49; $r++) {
    $realB += , ($ws.Cells.Item($r, 2).Value())
}
for ($r = 1; $r -le 49; $r++) {
    $ws.Cells.Item($r, 2).Value = (50 - $r)
}

$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add2($ws.Range("B1:B49"), 0, 2, 0, 0) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:D49")) | Out-Null
$ws.Sort.Header = -4142
$ws.Sort.Apply() | Out-Null

for ($r = 1; $r -le 49; $r++) {
    $ws.Cells.Item($r, 2).Value = $realB[$r-1]
}

# Restore the prior selection/scroll focus to reflect where the edit happened.
$ws.Range("E43").Select() | Out-Null
